# Update Bus_Name values in the two comparison sheets.
$wb = $excel.ActiveWorkbook

# --- DC_Comparison sheet: rows 2-4, column B ---
$wsDC = $wb.Worksheets.Item("DC_Comparison")
$dcNames = @("dcBus_3", "dcBus_2", "dcBus_1")
for ($i = 0; $i -lt $dcNames.Length; $i++) {
    $wsDC.Cells.Item($i + 2, 2).Value = $dcNames[$i]
}

# --- AC_Comparison sheet: rows 2-71, column B ---
$wsAC = $wb.Worksheets.Item("AC_Comparison")
$acNames = @(
    "Bus_33","Bus_56","Bus_38","Bus_37","Bus_54","Bus_17","Bus_18","Bus_6","Bus_42","Bus_48",
    "Bus_67","Bus_46","Bus_41","Bus_13","Bus_40","Bus_31","Bus_62","Bus_63","Bus_66","Bus_68",
    "Bus_58","Bus_57","Bus_55","Bus_52","Bus_51","Bus_47","Bus_45","Bus_15","Bus_19","Bus_24",
    "Bus_23","Bus_25","Bus_27","Bus_29","Bus_32","Bus_28","Bus_22","Bus_21","Bus_30","Bus_20",
    "Bus_26","Bus_16","Bus_70","Bus_2","Bus_9","Bus_10","Bus_4","Bus_5","Bus_11","Bus_12",
    "Bus_14","Bus_59","Bus_50","Bus_49","Bus_44","Bus_43","Bus_69","Bus_64","Bus_39","Bus_60",
    "Bus_65","Bus_8","Bus_1","Bus_3","Bus_7","Bus_53","Bus_61","Bus_35","Bus_36","Bus_34"
)
for ($i = 0; $i -lt $acNames.Length; $i++) {
    $wsAC.Cells.Item($i + 2, 2).Value = $acNames[$i]
}
